$d = $word.ActiveDocument

# 1. Update all "Updated 2023-10-29" footer lines to "Updated 2023-10-31"
$find = $d.Content.Find
$find.Execute("Updated 2023-10-29", $false, $false, $false, $false, $false, $true, 1, $false, "Updated 2023-10-31", 2)

# 2. Append a new "Insurance?" paragraph (BodyText style) after the last table
$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertParagraphAfter()

$insuranceXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Insurance?</w:t></w:r></w:p>'
$insurancePara = $d.Paragraphs.Last
$insurancePara.Range.InsertXML($insuranceXml)

# 3. Append the new Insurance table right after the "Insurance?" paragraph
$insurancePara2 = $d.Paragraphs.Last
$insurancePara2.Range.InsertParagraphAfter()

$tableXml = @'
    <w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:tblPr>
        <w:tblCellMar>
          <w:top w:w="0" w:type="dxa"/>
          <w:bottom w:w="0" w:type="dxa"/>
          <w:start w:w="60" w:type="dxa"/>
          <w:end w:w="60" w:type="dxa"/>
        </w:tblCellMar>
        <w:tblW w:type="auto" w:w="0"/>
        <w:tblLook w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:noHBand="0" w:noVBand="0"/>
        <w:jc w:val="center"/>
      </w:tblPr>
      <w:tr>
        <w:trPr>
          <w:cantSplit/>
          <w:tblHeader/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcBorders>
              <w:top w:val="single" w:sz="16" w:space="0" w:color="D3D3D3"/>
              <w:bottom w:val="single" w:sz="16" w:space="0" w:color="D3D3D3"/>
              <w:start w:val="single" w:space="0" w:color="D3D3D3"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="60"/>
              <w:keepNext/>
              <w:jc w:val="start"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="default">Response</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcBorders>
              <w:top w:val="single" w:sz="16" w:space="0" w:color="D3D3D3"/>
              <w:bottom w:val="single" w:sz="16" w:space="0" w:color="D3D3D3"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="60"/>
              <w:keepNext/>
              <w:jc w:val="end"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="default">Count</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcBorders>
              <w:top w:val="single" w:sz="16" w:space="0" w:color="D3D3D3"/>
              <w:bottom w:val="single" w:sz="16" w:space="0" w:color="D3D3D3"/>
              <w:end w:val="single" w:space="0" w:color="D3D3D3"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="60"/>
              <w:keepNext/>
              <w:jc w:val="end"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="default">Percentage</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:cantSplit/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcBorders>
              <w:top w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:bottom w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:start w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:end w:val="single" w:space="0" w:color="D3D3D3"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="60"/>
              <w:keepNext/>
              <w:jc w:val="start"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="default">Commercial Insurance</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcBorders>
              <w:top w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:bottom w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:start w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:end w:val="single" w:space="0" w:color="D3D3D3"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="60"/>
              <w:keepNext/>
              <w:jc w:val="end"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="default">60</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcBorders>
              <w:top w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:bottom w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:start w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:end w:val="single" w:space="0" w:color="D3D3D3"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="60"/>
              <w:keepNext/>
              <w:jc w:val="end"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="default">72.29</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:cantSplit/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcBorders>
              <w:top w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:bottom w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:start w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:end w:val="single" w:space="0" w:color="D3D3D3"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="60"/>
              <w:keepNext/>
              <w:jc w:val="start"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="default">Medicare</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcBorders>
              <w:top w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:bottom w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:start w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:end w:val="single" w:space="0" w:color="D3D3D3"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="60"/>
              <w:keepNext/>
              <w:jc w:val="end"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="default">7</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcBorders>
              <w:top w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:bottom w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:start w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:end w:val="single" w:space="0" w:color="D3D3D3"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="60"/>
              <w:keepNext/>
              <w:jc w:val="end"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="default">8.43</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:cantSplit/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcBorders>
              <w:top w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:bottom w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:start w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:end w:val="single" w:space="0" w:color="D3D3D3"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="60"/>
              <w:keepNext/>
              <w:jc w:val="start"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="default">Medicare Advantage</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcBorders>
              <w:top w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:bottom w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:start w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:end w:val="single" w:space="0" w:color="D3D3D3"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="60"/>
              <w:keepNext/>
              <w:jc w:val="end"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="default">4</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcBorders>
              <w:top w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:bottom w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:start w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:end w:val="single" w:space="0" w:color="D3D3D3"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="60"/>
              <w:keepNext/>
              <w:jc w:val="end"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="default">4.82</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:cantSplit/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcBorders>
              <w:top w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:bottom w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:start w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:end w:val="single" w:space="0" w:color="D3D3D3"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="60"/>
              <w:keepNext/>
              <w:jc w:val="start"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="default">Medicaid</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcBorders>
              <w:top w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:bottom w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:start w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:end w:val="single" w:space="0" w:color="D3D3D3"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="60"/>
              <w:keepNext/>
              <w:jc w:val="end"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="default">7</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcBorders>
              <w:top w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:bottom w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:start w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:end w:val="single" w:space="0" w:color="D3D3D3"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="60"/>
              <w:keepNext/>
              <w:jc w:val="end"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="default">8.43</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:cantSplit/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcBorders>
              <w:top w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:bottom w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:start w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:end w:val="single" w:space="0" w:color="D3D3D3"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="60"/>
              <w:keepNext/>
              <w:jc w:val="start"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="default">Self-Pay (No Insurance)</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcBorders>
              <w:top w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:bottom w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:start w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:end w:val="single" w:space="0" w:color="D3D3D3"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="60"/>
              <w:keepNext/>
              <w:jc w:val="end"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="default">5</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcBorders>
              <w:top w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:bottom w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:start w:val="single" w:space="0" w:color="D3D3D3"/>
              <w:end w:val="single" w:space="0" w:color="D3D3D3"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="60"/>
              <w:keepNext/>
              <w:jc w:val="end"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="default">6.02</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:cantSplit/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:gridSpan w:val="3"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:after="60"/>
              <w:keepNext/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
                <w:sz w:val="20"/>
              </w:rPr>
              <w:t xml:space="default">Updated 2023-10-31</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
'@

$tableTarget = $d.Paragraphs.Last
$tableTarget.Range.InsertXML($tableXml)

Write-Output ("Tables count: " + $d.Tables.Count)
Write-Output ("Paragraphs count: " + $d.Paragraphs.Count)
